$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the formatting of the
# existing header row (bold, centered, bordered) by copying H1's format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I ("I0") and J ("IF") - same value repeated in both
# columns for each data row.
$values = @{
    2  = 8
    3  = 8
    4  = 6
    5  = 7
    6  = 6
    7  = 6
    8  = 9
    9  = 8
    10 = 5
    11 = 4
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
